$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F15").Value = 38
$ws.Range("G15").Value = 14750.46
$ws.Range("F17").Value = 46
$ws.Range("G17").Value = 16738.48
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 4942.41
$ws.Range("B25").Value = 80690.17
$ws.Range("F69").Value = 312
$ws.Range("G69").Value = 35000.16
$ws.Range("F75").Value = 35
$ws.Range("G75").Value = 2461.2
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("F90").Value = 38
$ws.Range("G90").Value = 4303.88
$ws.Range("F92").Value = 79
$ws.Range("G92").Value = 3531.3
$ws.Range("B95").Value = 128689.17
$ws.Range("F126").Value = 75
$ws.Range("G126").Value = 3711
$ws.Range("F131").Value = 31
$ws.Range("G131").Value = 1311.92
$ws.Range("B135").Value = 17865.85
$ws.Range("F147").Value = 22
$ws.Range("G147").Value = 3386.9
$ws.Range("F153").Value = 40
$ws.Range("G153").Value = 6737.6
$ws.Range("F154").Value = 78
$ws.Range("G154").Value = 3957.72
$ws.Range("B168").Value = 32057.67
$ws.Range("F201").Value = 35
$ws.Range("G201").Value = 3570.35
$ws.Range("F202").Value = 176
$ws.Range("G202").Value = 12622.72
$ws.Range("F203").Value = 51
$ws.Range("G203").Value = 4543.59
$ws.Range("F207").Value = 118
$ws.Range("G207").Value = 8767.4
$ws.Range("B208").Value = 34470.8
$ws.Range("F210").Value = 169
$ws.Range("G210").Value = 19759.48
$ws.Range("F211").Value = 1408
$ws.Range("G211").Value = 26048
$ws.Range("F217").Value = 12
$ws.Range("G217").Value = 1675.68
$ws.Range("B222").Value = 53934.31
$ws.Range("F239").Value = 0
$ws.Range("G239").Value = 0
$ws.Range("F240").Value = 35
$ws.Range("G240").Value = 1093.75
$ws.Range("F250").Value = 5
$ws.Range("G250").Value = 234.35
$ws.Range("F255").Value = 14
$ws.Range("G255").Value = 7770.42
$ws.Range("B258").Value = 43986.14
$ws.Range("F289").Value = 93
$ws.Range("G289").Value = 2457.06
$ws.Range("B290").Value = 70570.33
$ws.Range("F305").Value = 40
$ws.Range("G305").Value = 1075.6
$ws.Range("B307").Value = 7140.52
$ws.Range("F312").Value = 115
$ws.Range("G312").Value = 6951.75
$ws.Range("F313").Value = 75
$ws.Range("G313").Value = 1911
$ws.Range("F317").Value = 77
$ws.Range("G317").Value = 4322.01
$ws.Range("F328").Value = 43
$ws.Range("G328").Value = 2038.2
$ws.Range("F330").Value = 18
$ws.Range("G330").Value = 853.2
$ws.Range("F332").Value = 47
$ws.Range("G332").Value = 1894.1
$ws.Range("F336").Value = 67
$ws.Range("G336").Value = 12491.48
$ws.Range("F340").Value = 103
$ws.Range("G340").Value = 2964.34
$ws.Range("F342").Value = 48
$ws.Range("G342").Value = 3537.12
$ws.Range("F346").Value = 47
$ws.Range("G346").Value = 1686.83
$ws.Range("B347").Value = 121834.59
$ws.Range("F383").Value = 12
$ws.Range("G383").Value = 497.04
$ws.Range("F386").Value = 116
$ws.Range("G386").Value = 1095.04
$ws.Range("F387").Value = 72
$ws.Range("G387").Value = 1585.44
$ws.Range("B393").Value = 8924.440000000001
$ws.Range("F423").Value = 17
$ws.Range("G423").Value = 3020.22
$ws.Range("B425").Value = 3088.75
$ws.Range("F439").Value = 39
$ws.Range("G439").Value = 2140.32
$ws.Range("B455").Value = 43567.18
$ws.Range("F461").Value = 78
$ws.Range("G461").Value = 2582.58
$ws.Range("F462").Value = 264
$ws.Range("G462").Value = 4187.04
$ws.Range("F463").Value = 58
$ws.Range("G463").Value = 1920.38
$ws.Range("B470").Value = 10260.57
$ws.Range("B607").Value = 1581135.81
$ws.Range("B608").Value = 1581135.81
